$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 555 (shifts existing rows 555..581 down to 556..582)
$ws.Rows.Item(555).Insert()

# Populate the newly inserted row 555 with the new record
$ws.Range("A555").Value = 11
$ws.Range("B555").Value = "Vega Monumental Concepción"
$ws.Range("C555").Value = "Bíobío"
$ws.Range("D555").Value = 45267
$ws.Range("E555").Value = 8
$ws.Range("F555").Value = "Fruta"
$ws.Range("G555").Value = 100102
$ws.Range("H555").Value = "Cítricos"
$ws.Range("I555").Value = 100102005
$ws.Range("J555").Value = "Naranja"
$ws.Range("K555").Value = "Valencia"
$ws.Range("L555").Value = "Primera"
$ws.Range("M555").Value = 350
$ws.Range("N555").Value = 12000
$ws.Range("O555").Value = 12000
$ws.Range("P555").Value = 12000
$ws.Range("Q555").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R555").Value = "Región de O'Higgins"
$ws.Range("S555").Value = 800
$ws.Range("T555").Value = 15

# Make sure the date cell keeps the date number format used by the rest of column D
$ws.Range("D555").NumberFormat = $ws.Range("D554").NumberFormat
